# Update calculated error values on "Finite wing theory" and "Slender wing theory"
# sheets (FlatWing results that were forgotten to be pushed).

$wb = $excel.ActiveWorkbook

$wsFinite = $wb.Worksheets.Item("Finite wing theory")
$wsFinite.Range("D2").Value = 0.38013691964065927
$wsFinite.Range("E2").Value = 0.54737570206066166
$wsFinite.Range("F2").Value = 0.31982659424771431

$wsFinite.Range("D3").Value = 0.064075224142106035
$wsFinite.Range("E3").Value = 0.53561436446007782
$wsFinite.Range("F3").Value = 0.19758996236955251

$wsFinite.Range("D4").Value = 0.064511134667174883
$wsFinite.Range("E4").Value = 0.36125841268767245
$wsFinite.Range("F4").Value = 0.13298607494057224

$wsFinite.Range("D5").Value = 0.13888180900026556
$wsFinite.Range("E5").Value = 0.12183814259614477
$wsFinite.Range("F5").Value = 0.095159108012391469

$wsSlender = $wb.Worksheets.Item("Slender wing theory")
$wsSlender.Range("D2").Value = 0.055864314405044503
$wsSlender.Range("E2").Value = 0.78818225375800222
$wsSlender.Range("F2").Value = 0.0095273377258247027

$wsSlender.Range("D3").Value = 0.31662228116045604
$wsSlender.Range("E3").Value = 0.8084616487752142
$wsSlender.Range("F3").Value = 0.029944907984831848

$wsSlender.Range("D4").Value = 0.53508529616933442
$wsSlender.Range("E4").Value = 0.87816520415989308
$wsSlender.Range("F4").Value = 0.050806579967106746

$wsSlender.Range("D5").Value = 0.67901447670584369
$wsSlender.Range("E5").Value = 0.93024294310395828
$wsSlender.Range("F5").Value = 0.063021538037665248

$wb.Save()
